$d = $word.ActiveDocument

# Locate the "3 " text that holds the "User story" number near the top of
# the document (e.g. "User story 3   Thème : ...").
$findRange = $d.Content
$found = $findRange.Find.Execute("3 ", $true, $false, $false, $false, $false, $true)
if (-not $found) {
    throw "Could not locate the 'User story' number text ('3 ') to update."
}

$idx = $findRange.Start

# Position right after the "3" character itself (between "3" and the space
# that follows it) - this is where the cursor ends up once "3" becomes "4".
$bmPos = $idx + 1

# Drop the (hidden) "_GoBack" bookmark at that new edit location, mirroring
# how Word marks the last place you edited. A bookmark named "_GoBack" can
# only exist once in the document, so re-adding it here automatically moves
# it away from its previous location (next to "(de la dernière heure)").
$rBookmark = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $rBookmark) | Out-Null

# Finally, change the "3" character itself to "4" (the space stays as-is).
$rDigit = $d.Range($idx, $idx + 1)
$rDigit.Text = "4"
